$d = $word.ActiveDocument

$pairs = @(
    @("97×29=", "20×22="),
    @("94×73=", "49×54="),
    @("48×94=", "86×93="),
    @("81×12=", "42×15="),
    @("15×69=", "51×98="),
    @("93×56=", "74×78="),
    @("31×23=", "59×42="),
    @("51×51=", "86×33="),
    @("95×37=", "13×35="),
    @("48×35=", "29×56="),
    @("93×18=", "85×68="),
    @("45×33=", "57×35="),
    @("76×66=", "58×42="),
    @("88×84=", "62×83="),
    @("66×47=", "74×87="),
    @("32×48=", "21×62="),
    @("68×93=", "35×50="),
    @("40×14=", "73×24="),
    @("27×11=", "83×77="),
    @("34×29=", "33×27="),
    @("67×12=", "60×84="),
    @("61×58=", "20×51="),
    @("98×47=", "48×24="),
    @("34×95=", "75×50="),
    @("30×53=", "25×19="),
    @("29×79=", "98×23="),
    @("60×17=", "71×98="),
    @("67×70=", "78×12="),
    @("68×57=", "66×16="),
    @("41×72=", "79×26="),
    @("55×69=", "75×47="),
    @("82×76=", "55×36="),
    @("58×61=", "54×31="),
    @("23×98=", "92×79="),
    @("82×18=", "32×95="),
    @("68×79=", "61×53="),
    @("49×48=", "44×91="),
    @("60×76=", "66×68="),
    @("24×28=", "67×83="),
    @("44×22=", "21×23="),
    @("56×76=", "72×74="),
    @("10×54=", "87×71="),
    @("91×83=", "94×52="),
    @("54×84=", "67×44="),
    @("100×20=", "57×92="),
    @("11×65=", "66×38="),
    @("18×71=", "100×12="),
    @("22×15=", "15×89="),
    @("18×79=", "41×10="),
    @("25×17=", "45×43="),
    @("79×51=", "58×56="),
    @("58×18=", "99×70="),
    @("100×48=", "75×38="),
    @("60×35=", "88×96="),
    @("35×18=", "11×94="),
    @("80×82=", "83×50="),
    @("16×43=", "90×16="),
    @("64×94=", "97×19="),
    @("66×54=", "19×52="),
    @("95×57=", "13×98="),
    @("97×42=", "60×52="),
    @("79×11=", "20×97="),
    @("16×75=", "40×28="),
    @("97×41=", "73×71="),
    @("43×65=", "46×22="),
    @("25×60=", "98×93="),
    @("45×62=", "58×22="),
    @("96×36=", "58×93="),
    @("89×19=", "16×84="),
    @("79×36=", "73×90="),
    @("33×32=", "55×30="),
    @("85×98=", "42×63="),
    @("25×33=", "11×29="),
    @("61×79=", "69×30="),
    @("90×11=", "81×39="),
    @("52×74=", "63×50="),
    @("72×76=", "18×100="),
    @("25×82=", "100×32="),
    @("70×60=", "85×15="),
    @("71×34=", "85×41="),
    @("22×86=", "32×20="),
    @("70×40=", "97×63="),
    @("62×13=", "55×62="),
    @("52×68=", "51×12="),
    @("89×17=", "11×64="),
    @("61×51=", "95×55="),
    @("92×62=", "58×83="),
    @("27×73=", "51×20="),
    @("99×57=", "92×32="),
    @("100×21=", "25×21="),
    @("83×37=", "82×29="),
    @("43×28=", "39×39="),
    @("91×53=", "48×53="),
    @("66×43=", "96×75="),
    @("38×49=", "70×88="),
    @("33×64=", "37×97="),
    @("66×60=", "18×53="),
    @("77×76=", "11×16="),
    @("58×96=", "61×41="),
    @("70×32=", "32×80="),
)

$replaced = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($ok) { $replaced++ }
}
Write-Host "Replaced $replaced of $($pairs.Count) expressions"